# "Final changes before McGill Sim"
#
# 1) Add an "Occupancy Detection" value of "PIR" for every lab row
#    (N4:N36) on the "laboratories" sheet.
# 2) Make "laboratories" the active/selected sheet (it was "fumehoods"
#    before), and leave the cursor on N9 there.

$wb = $excel.ActiveWorkbook

$wsLab = $wb.Worksheets.Item("laboratories")

# --- 1. Fill in the new "Occupancy Detection" column -----------------------
for ($row = 4; $row -le 36; $row++) {
    $wsLab.Range("N$row").Value = "PIR"
}

# --- 2. Switch the active tab to "laboratories" and update its selection ---
$wsLab.Activate() | Out-Null
$wsLab.Range("N9").Select() | Out-Null
